$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 5345.4546  # H113: was 4799.923
$ws.Cells.Item(113, 9).Value = 4880  # I113: was 4889.9
$ws.Cells.Item(113, 10).Value = 10000  # J113: was 4500
$ws.Cells.Item(113, 11).Value = 4880  # K113: was 4889.9
$ws.Cells.Item(113, 12).Value = 10000  # L113: was 4500
$ws.Cells.Item(113, 13).Value = -1626  # M113: was -1635.9
$ws.Cells.Item(113, 14).Value = -16508  # N113: was -11008
$ws.Cells.Item(128, 8).Value = 40000  # H128: was 29655.555
$ws.Cells.Item(128, 10).Value = 40000  # J128: was 29655.555
$ws.Cells.Item(128, 12).Value = 40000  # L128: was 29655.555
$ws.Cells.Item(128, 14).Value = -49960  # N128: was -39615.555
$ws.Cells.Item(132, 8).Value = 4168493.8  # H132: was 4257173.5
$ws.Cells.Item(132, 9).Value = 5129720  # I132: was 5264698
$ws.Cells.Item(132, 11).Value = 15389160  # K132: was 15794094
$ws.Cells.Item(132, 13).Value = -15386630  # M132: was -15791564
$ws.Cells.Item(137, 8).Value = 2706354.8  # H137: was 2567575
$ws.Cells.Item(137, 9).Value = 4004167.5  # I137: was 3229286.8
$ws.Cells.Item(137, 10).Value = 2578.1667  # J137: was 3442.875
$ws.Cells.Item(137, 11).Value = 12012502.5  # K137: was 9687860.399999999
$ws.Cells.Item(137, 12).Value = 7734.500100000001  # L137: was 10328.625
$ws.Cells.Item(137, 13).Value = -12009952.5  # M137: was -9685310.399999999
$ws.Cells.Item(137, 14).Value = -12834.5001  # N137: was -15428.625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 19233170  # H2: was 11906573
$ws.Cells.Item(2, 9).Value = 31251452  # I2: was 16667862
$ws.Cells.Item(2, 10).Value = 3920  # J2: was 3350
$ws.Cells.Item(2, 11).Value = 31251452  # K2: was 16667862
$ws.Cells.Item(2, 12).Value = 3920  # L2: was 3350
$ws.Cells.Item(2, 13).Value = -31251339  # M2: was -16667749
$ws.Cells.Item(2, 14).Value = -4146  # N2: was -3576
$ws.Cells.Item(32, 8).Value = 6984.7534  # H32: was 6949.1
$ws.Cells.Item(32, 9).Value = 5819.411  # I32: was 5858.92
$ws.Cells.Item(32, 10).Value = 28252.25  # J32: was 23301.8
$ws.Cells.Item(32, 11).Value = 5819.411  # K32: was 5858.92
$ws.Cells.Item(32, 12).Value = 28252.25  # L32: was 23301.8
$ws.Cells.Item(32, 13).Value = -5532.411  # M32: was -5571.92
$ws.Cells.Item(32, 14).Value = -28826.25  # N32: was -23875.8
$ws.Cells.Item(61, 8).Value = 1915.9429  # H61: was 2120.625
$ws.Cells.Item(61, 9).Value = 798.13336  # I61: was 995
$ws.Cells.Item(61, 10).Value = 8622.799999999999  # J61: was 10000
$ws.Cells.Item(61, 11).Value = 798.13336  # K61: was 995
$ws.Cells.Item(61, 12).Value = 8622.799999999999  # L61: was 10000
$ws.Cells.Item(61, 13).Value = -586.13336  # M61: was -783
$ws.Cells.Item(61, 14).Value = -9046.799999999999  # N61: was -10424
$ws.Cells.Item(97, 8).Value = 675  # H97: was 709.7727
$ws.Cells.Item(97, 9).Value = 590.2381  # I97: was 606.75
$ws.Cells.Item(97, 10).Value = 1120  # J97: was 1740
$ws.Cells.Item(97, 11).Value = 590.2381  # K97: was 606.75
$ws.Cells.Item(97, 12).Value = 1120  # L97: was 1740
$ws.Cells.Item(97, 13).Value = -94.23810000000003  # M97: was -110.75
$ws.Cells.Item(97, 14).Value = -2112  # N97: was -2732
$ws.Cells.Item(116, 8).Value = 19233170  # H116: was 11906573
$ws.Cells.Item(116, 9).Value = 31251452  # I116: was 16667862
$ws.Cells.Item(116, 10).Value = 3920  # J116: was 3350
$ws.Cells.Item(116, 11).Value = 31251452  # K116: was 16667862
$ws.Cells.Item(116, 12).Value = 3920  # L116: was 3350
$ws.Cells.Item(116, 13).Value = -31249158  # M116: was -16665568
$ws.Cells.Item(116, 14).Value = -8508  # N116: was -7938
$ws.Cells.Item(136, 8).Value = 1915.9429  # H136: was 2120.625
$ws.Cells.Item(136, 9).Value = 798.13336  # I136: was 995
$ws.Cells.Item(136, 10).Value = 8622.799999999999  # J136: was 10000
$ws.Cells.Item(136, 11).Value = 2394.40008  # K136: was 2985
$ws.Cells.Item(136, 12).Value = 25868.4  # L136: was 30000
$ws.Cells.Item(136, 13).Value = 155.5999199999997  # M136: was -435
$ws.Cells.Item(136, 14).Value = -30968.4  # N136: was -35100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 19233170  # H3: was 11906573
$ws.Cells.Item(3, 9).Value = 31251452  # I3: was 16667862
$ws.Cells.Item(3, 10).Value = 3920  # J3: was 3350
$ws.Cells.Item(3, 11).Value = 31251452  # K3: was 16667862
$ws.Cells.Item(3, 12).Value = 3920  # L3: was 3350
$ws.Cells.Item(3, 13).Value = -31251338  # M3: was -16667748
$ws.Cells.Item(3, 14).Value = -4148  # N3: was -3578
$ws.Cells.Item(98, 8).Value = 0  # H98: was 28314
$ws.Cells.Item(98, 10).Value = 0  # J98: was 28314
$ws.Cells.Item(98, 12).Value = 0  # L98: was 28314
$ws.Cells.Item(98, 14).Value = $null  # N98: was -34304, deleted
$ws.Cells.Item(107, 8).Value = 2586.182  # H107: was 3023.111
$ws.Cells.Item(107, 9).Value = 816.875  # I107: was 882.5
$ws.Cells.Item(107, 11).Value = 816.875  # K107: was 882.5
$ws.Cells.Item(107, 13).Value = 1103.125  # M107: was 1037.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1614790.4  # H31: was 1787748.4
$ws.Cells.Item(31, 9).Value = 1819449.5  # I31: was 2042176.1
$ws.Cells.Item(31, 11).Value = 1819449.5  # K31: was 2042176.1
$ws.Cells.Item(31, 13).Value = -1819154.5  # M31: was -2041881.1
$ws.Cells.Item(34, 8).Value = 1614790.4  # H34: was 1787748.4
$ws.Cells.Item(34, 9).Value = 1819449.5  # I34: was 2042176.1
$ws.Cells.Item(34, 11).Value = 1819449.5  # K34: was 2042176.1
$ws.Cells.Item(34, 13).Value = -1819247.5  # M34: was -2041974.1
$ws.Cells.Item(58, 8).Value = 17245480  # H58: was 13516893
$ws.Cells.Item(58, 9).Value = 2800.95  # I58: was 2389.8
$ws.Cells.Item(58, 10).Value = 55562544  # J58: was 41672108
$ws.Cells.Item(58, 11).Value = 2800.95  # K58: was 2389.8
$ws.Cells.Item(58, 12).Value = 55562544  # L58: was 41672108
$ws.Cells.Item(58, 13).Value = -2597.95  # M58: was -2186.8
$ws.Cells.Item(58, 14).Value = -55562950  # N58: was -41672514
$ws.Cells.Item(99, 8).Value = 3420.1428  # H99: was 3040.818
$ws.Cells.Item(99, 9).Value = 1989.2727  # I99: was 1888.625
$ws.Cells.Item(99, 10).Value = 8666.666999999999  # J99: was 6113.3335
$ws.Cells.Item(99, 11).Value = 1989.2727  # K99: was 1888.625
$ws.Cells.Item(99, 12).Value = 8666.666999999999  # L99: was 6113.3335
$ws.Cells.Item(99, 13).Value = -491.2727  # M99: was -390.625
$ws.Cells.Item(99, 14).Value = -11662.667  # N99: was -9109.333500000001
$ws.Cells.Item(107, 8).Value = 1717.6875  # H107: was 1845.6
$ws.Cells.Item(107, 9).Value = 543.9231  # I107: was 564.25
$ws.Cells.Item(107, 10).Value = 6804  # J107: was 6971
$ws.Cells.Item(107, 11).Value = 543.9231  # K107: was 564.25
$ws.Cells.Item(107, 12).Value = 6804  # L107: was 6971
$ws.Cells.Item(107, 13).Value = 1376.0769  # M107: was 1355.75
$ws.Cells.Item(107, 14).Value = -10644  # N107: was -10811
$ws.Cells.Item(126, 8).Value = 3420.1428  # H126: was 3040.818
$ws.Cells.Item(126, 9).Value = 1989.2727  # I126: was 1888.625
$ws.Cells.Item(126, 10).Value = 8666.666999999999  # J126: was 6113.3335
$ws.Cells.Item(126, 11).Value = 5967.8181  # K126: was 5665.875
$ws.Cells.Item(126, 12).Value = 26000.001  # L126: was 18340.0005
$ws.Cells.Item(126, 13).Value = -3497.8181  # M126: was -3195.875
$ws.Cells.Item(126, 14).Value = -30940.001  # N126: was -23280.0005
$ws.Cells.Item(132, 8).Value = 2034.7059  # H132: was 2221.8462
$ws.Cells.Item(132, 9).Value = 1683.2258  # I132: was 1759
$ws.Cells.Item(132, 10).Value = 5666.6665  # J132: was 3400
$ws.Cells.Item(132, 11).Value = 5049.6774  # K132: was 5277
$ws.Cells.Item(132, 12).Value = 16999.9995  # L132: was 10200
$ws.Cells.Item(132, 13).Value = -2519.6774  # M132: was -2747
$ws.Cells.Item(132, 14).Value = -22059.9995  # N132: was -15260
$ws.Cells.Item(136, 8).Value = 17245480  # H136: was 13516893
$ws.Cells.Item(136, 9).Value = 2800.95  # I136: was 2389.8
$ws.Cells.Item(136, 10).Value = 55562544  # J136: was 41672108
$ws.Cells.Item(136, 11).Value = 8402.849999999999  # K136: was 7169.400000000001
$ws.Cells.Item(136, 12).Value = 166687632  # L136: was 125016324
$ws.Cells.Item(136, 13).Value = -5852.849999999999  # M136: was -4619.400000000001
$ws.Cells.Item(136, 14).Value = -166692732  # N136: was -125021424
$ws.Cells.Item(141, 8).Value = 13384.211  # H141: was 16300
$ws.Cells.Item(141, 10).Value = 13384.211  # J141: was 16300
$ws.Cells.Item(141, 12).Value = 13384.211  # L141: was 16300
$ws.Cells.Item(141, 14).Value = -23744.211  # N141: was -26660

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 1257.8235  # H122: was 1339.8125
$ws.Cells.Item(122, 9).Value = 598.3333  # I122: was 679.875
$ws.Cells.Item(122, 11).Value = 5384.9997  # K122: was 6118.875
$ws.Cells.Item(122, 13).Value = -2934.9997  # M122: was -3668.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1889.5454  # H113: was 1931.9048
$ws.Cells.Item(113, 10).Value = 2857.1428  # J113: was 3166.6667
$ws.Cells.Item(113, 12).Value = 2857.1428  # L113: was 3166.6667
$ws.Cells.Item(113, 14).Value = -7197.1428  # N113: was -7506.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 80007  # H15: was 28480
$ws.Cells.Item(15, 9).Value = 0  # I15: was 4006
$ws.Cells.Item(15, 10).Value = 80007  # J15: was 36638
$ws.Cells.Item(15, 11).Value = 0  # K15: was 4006
$ws.Cells.Item(15, 12).Value = 80007  # L15: was 36638
$ws.Cells.Item(15, 13).Value = $null  # M15: was -3718, deleted
$ws.Cells.Item(15, 14).Value = -80583  # N15: was -37214
$ws.Cells.Item(107, 8).Value = 1540  # H107: was 710.0417
$ws.Cells.Item(107, 9).Value = 512.5  # I107: was 301.2353
$ws.Cells.Item(107, 10).Value = 5650  # J107: was 1702.8572
$ws.Cells.Item(107, 11).Value = 1537.5  # K107: was 903.7058999999999
$ws.Cells.Item(107, 12).Value = 16950  # L107: was 5108.571599999999
$ws.Cells.Item(107, 13).Value = 382.5  # M107: was 1016.2941
$ws.Cells.Item(107, 14).Value = -20790  # N107: was -8948.571599999999
$ws.Cells.Item(132, 8).Value = 152718.38  # H132: was 162604.97
$ws.Cells.Item(132, 9).Value = 176094.9  # I132: was 189424.62
$ws.Cells.Item(132, 11).Value = 528284.7  # K132: was 568273.86
$ws.Cells.Item(132, 13).Value = -525754.7  # M132: was -565743.86
